# msz - Login-Testcase success for mobile + web
#
# Row 1 (headers) stays the same content, rows 2-3 get new
# "BaseState..." data instead of the old Variable1/Variable2 placeholder
# rows, and a new row 4 (AUT / Chromium / Pixel9Pro_API35) is appended.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: was Variable1 / Wert1 / Value1
$ws.Range("A2").Value = "BaseStateChromium"
$ws.Range("B2").Value = "Started"
$ws.Range("C2").Value = "Not Started"

# Row 3: was Variable2 / Wert2 / Value2
$ws.Range("A3").Value = "BaseStatePixel9Pro_API35"
$ws.Range("B3").Value = "Not Started"
$ws.Range("C3").Value = "Started"

# Row 4: new row for the AUT / browser / device combo
$ws.Range("A4").Value = "AUT"
$ws.Range("B4").Value = "Chromium"
$ws.Range("C4").Value = "Pixel9Pro_API35"

# Widen columns B and C to fit the new content (best effort - matches the
# author's manual resize / best-fit of columns B and C).
$ws.Columns.Item(2).ColumnWidth = 16
$ws.Columns.Item(3).ColumnWidth = 13.25

# Selection moves to the newly added row's first cell.
$ws.Range("A4").Select() | Out-Null
